$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "'1120194100370"
$ws.Range("C2").Value = 100

$ws.Range("B3").Value = "'1120170200907"
$ws.Range("C3").Value = 100

$ws.Range("G7").Select()
